$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = -1
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 2
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = -4
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = 3
